$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    $header = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>'
    $footer = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $header + $bodyXml + $footer
}

# Paragraph 1: "This is a test for the docx4j," - drop the pPr (paragraph-mark rFonts hint)
$p1Xml = '<w:p><w:r><w:t>T</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>his is a test for the docx4j</w:t></w:r><w:r w:rsidR="005C1AFD"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>,</w:t></w:r></w:p>'
$d.Paragraphs.Item(1).Range.InsertXML((New-PkgXml $p1Xml)) | Out-Null

# Paragraph 2: "You name is  gaobin  ," -> drop pPr, add a trailing space run after
# "ou name is", move the _GoBack bookmark here, and rename gaobin -> namevalue
$p2Xml = '<w:p><w:r><w:t>Y</w:t></w:r><w:r w:rsidR="0074615A"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>ou name is</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:name="_GoBack" w:id="0"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="009A5830"><w:rPr><w:rFonts w:hAnsi="Verdana" w:ascii="Verdana" w:hint="eastAsia"/><w:color w:val="333333"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>  namevalue  </w:t></w:r><w:r w:rsidR="009A5830"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="009A5830"><w:rPr><w:rFonts w:hAnsi="Verdana" w:ascii="Verdana" w:hint="eastAsia"/><w:color w:val="333333"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="005C1AFD"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>,</w:t></w:r></w:p>'
$d.Paragraphs.Item(2).Range.InsertXML((New-PkgXml $p2Xml)) | Out-Null

# Paragraph 3: "You address is  China  ," -> drop pPr, rename China -> addressvalue
$p3Xml = '<w:p><w:r><w:t>Y</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">ou address is </w:t></w:r><w:r w:rsidR="009A5830"><w:rPr><w:rFonts w:hAnsi="Verdana" w:ascii="Verdana" w:hint="eastAsia"/><w:color w:val="333333"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>  addressvalue  </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="009A5830"><w:rPr><w:rFonts w:hAnsi="Verdana" w:ascii="Verdana" w:hint="eastAsia"/><w:color w:val="333333"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="005C1AFD"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>,</w:t></w:r></w:p>'
$d.Paragraphs.Item(3).Range.InsertXML((New-PkgXml $p3Xml)) | Out-Null

# Paragraph 4: "That's all for this file. " -> drop the _GoBack bookmark (it moved to paragraph 2)
$p4Xml = '<w:p><w:r><w:t>T</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>hat</w:t></w:r><w:r><w:t>' + [char]0x2019 + '</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>s all for this file.</w:t></w:r><w:r w:rsidR="009A5830" w:rsidRPr="009A5830"><w:rPr><w:rFonts w:hAnsi="Verdana" w:ascii="Verdana" w:hint="eastAsia"/><w:color w:val="333333"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$d.Paragraphs.Item(4).Range.InsertXML((New-PkgXml $p4Xml)) | Out-Null

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
